$d = $word.ActiveDocument

$items = @(
    "Students & Administrators should be able to login using CAS",
    "Student should be able to access a dashboard with remaining/completed requirements",
    "Students should be able to upload supporting documentation",
    "Administrators should be able to approve uploaded documents",
    "Administrators should be able to look up a specific student to see completed/remaining requirements",
    "Administrators should be able to generate a report of all active STIC applicants and their completed/missing requirements.",
    "Student should be able to apply for the STIC program (preliminary approval)",
    "Administrators should be able to approve applicants",
    "Administrators should be able to access uploaded documents for review",
    "Administrators should be able to add/remove/edit requirements sets for degrees/programs"
)

# Grab the bullet-list template in use for the existing "Include new risks"
# item (numId 1) so the new bullets below can be appended to that same
# list rather than minting a brand-new list.
$lastPara = $d.Paragraphs.Last
$tmpl = $lastPara.Range.ListFormat.ListTemplate

# Add a plain (non-list) paragraph for the "Functional Requirements:"
# heading right after the last existing paragraph.
$r = $lastPara.Range.Duplicate
$r.Collapse(0)
$r.InsertParagraphAfter()
$headingPara = $d.Paragraphs.Last
$headingPara.Style = "Normal"
$headingPara.Range.Text = "Functional Requirements:"

# Append each functional-requirement item as a bulleted ListParagraph tied
# to the same list (numId 1, level 0) as the rest of the document's bullets.
foreach ($item in $items) {
    $prevPara = $d.Paragraphs.Last
    $r = $prevPara.Range.Duplicate
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $p = $d.Paragraphs.Last
    $p.Style = "List Paragraph"
    $p.Range.ListFormat.ApplyListTemplateWithLevel($tmpl, $true, 1, $false, 1)
    $p.Range.Text = $item
}

Write-Output "done"
